# Update predicted-price based compared-returns columns (G, H) for rows 2-57,
# and the first-row mean_return_pct_change (I2), to reflect the new
# prediction run referenced in the commit message
# ("new predicted prices and compared returns").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.04599563196556224
$ws.Range("H2").Value = 38.14585772397647
$ws.Range("I2").Value = -38.16679062961737

$ws.Range("G3").Value = 0.05815752889129102
$ws.Range("H3").Value = 15.24398503988901

$ws.Range("G4").Value = 0.04308162222390904
$ws.Range("H4").Value = -11.49276895599284

$ws.Range("G5").Value = 0.002850564834563522
$ws.Range("H5").Value = -95.96548980216572

$ws.Range("G6").Value = -0.1225407985429096
$ws.Range("H6").Value = -4.178611042672927

$ws.Range("G7").Value = -0.0851732567400029
$ws.Range("H7").Value = 31.87963469391404

$ws.Range("G8").Value = -0.2387438759557693
$ws.Range("H8").Value = -19.8641459726214

$ws.Range("G9").Value = -0.3046998927184322
$ws.Range("H9").Value = -0.868323808648648

$ws.Range("G10").Value = -0.01672526387665524
$ws.Range("H10").Value = -1414.140510368105

$ws.Range("G11").Value = 0.01831867858298696
$ws.Range("H11").Value = 189.0839665977422

$ws.Range("G12").Value = 0.1900226489853923
$ws.Range("H12").Value = -10.32914761763166

$ws.Range("G13").Value = 0.2325644630009688
$ws.Range("H13").Value = -1.034962421647844

$ws.Range("G14").Value = -0.09057546250880998
$ws.Range("H14").Value = 0.5214776706301707

$ws.Range("G15").Value = -0.08512714613283107
$ws.Range("H15").Value = -19.95021407457609

$ws.Range("G16").Value = 0.1887894375790652
$ws.Range("H16").Value = -1.391379320928736

$ws.Range("G17").Value = 0.187350766286823
$ws.Range("H17").Value = 7.753593027255913

$ws.Range("G18").Value = 0.03891736041782435
$ws.Range("H18").Value = -28.23628328627338

$ws.Range("G19").Value = 0.05301573480766496
$ws.Range("H19").Value = -38.38217052383815

$ws.Range("G20").Value = -0.03112155857827627
$ws.Range("H20").Value = -344.4951232309485

$ws.Range("G21").Value = -0.02396442743188552
$ws.Range("H21").Value = 55.51746458381414

$ws.Range("G22").Value = 0.07403760399341894
$ws.Range("H22").Value = 13.42806189613998

$ws.Range("G23").Value = 0.0903595973489226
$ws.Range("H23").Value = 56.67467940789406

$ws.Range("G24").Value = 0.03767523095279327
$ws.Range("H24").Value = 16.30271870560653

$ws.Range("G25").Value = 0.05919356130235929
$ws.Range("H25").Value = 101.1188470387414

$ws.Range("G26").Value = 0.09369504484397112
$ws.Range("H26").Value = -17.29731205608081

$ws.Range("G27").Value = 0.1053630909197535
$ws.Range("H27").Value = 16.82650806215831

$ws.Range("G28").Value = 0.1285394614181637
$ws.Range("H28").Value = 9.425909870319098

$ws.Range("G29").Value = 0.1153349065302509
$ws.Range("H29").Value = -3.58843248912103

$ws.Range("G30").Value = 0.07881122313039023
$ws.Range("H30").Value = 17.22917845653225

$ws.Range("G31").Value = 0.08176148935002298
$ws.Range("H31").Value = 19.1551663070107

$ws.Range("G32").Value = 0.04765114863968132
$ws.Range("H32").Value = 9.125044539007588

$ws.Range("G33").Value = 0.05472601292925366
$ws.Range("H33").Value = 0.7133126299306605

$ws.Range("G34").Value = -0.009968748171535105
$ws.Range("H34").Value = 47.80382559483382

$ws.Range("G35").Value = 0.02283512539432225
$ws.Range("H35").Value = 63.66134756407135

$ws.Range("G36").Value = -0.01821855986303961
$ws.Range("H36").Value = -217.8514596856645

$ws.Range("G37").Value = 0.01050770562199627
$ws.Range("H37").Value = -16.09724313001396

$ws.Range("G38").Value = 0.04106811338735054
$ws.Range("H38").Value = -42.75646677270537

$ws.Range("G39").Value = 0.0524457206853457
$ws.Range("H39").Value = 21.78975545446547

$ws.Range("G40").Value = 0.06633427153731296
$ws.Range("H40").Value = 48.300062175929

$ws.Range("G41").Value = 0.09111490214836002
$ws.Range("H41").Value = 637.2160089312318

$ws.Range("G42").Value = 0.06007985842396252
$ws.Range("H42").Value = 14.92163832333614

$ws.Range("G43").Value = 0.0605315251399902
$ws.Range("H43").Value = 21.3059404650695

$ws.Range("G44").Value = 0.1008047878330977
$ws.Range("H44").Value = -23.48832998534596

$ws.Range("G45").Value = 0.1017472072531099
$ws.Range("H45").Value = -43.29315812890521

$ws.Range("G46").Value = -0.04177109995115467
$ws.Range("H46").Value = 4.914067468415594

$ws.Range("G47").Value = -0.04047319334807199
$ws.Range("H47").Value = -1445.033839110868

$ws.Range("G48").Value = -0.004948056114730009
$ws.Range("H48").Value = -134.1437849572903

$ws.Range("G49").Value = 0.01489094951378256
$ws.Range("H49").Value = 367.8797604433791

$ws.Range("G50").Value = 0.1255119757279282
$ws.Range("H50").Value = -12.18730544353951

$ws.Range("G51").Value = 0.1463968632333001
$ws.Range("H51").Value = 11.78301329157617

$ws.Range("G52").Value = 0.08434940901631864
$ws.Range("H52").Value = 36.15129468426858

$ws.Range("G53").Value = 0.07187163582054601
$ws.Range("H53").Value = 17.47257747293467

$ws.Range("G54").Value = -0.1475315042842161
$ws.Range("H54").Value = -65.19601335629204

$ws.Range("G55").Value = -0.1037232182198249
$ws.Range("H55").Value = 0.06107936907795747

$ws.Range("G56").Value = 0.1257906606580994
$ws.Range("H56").Value = -18.9094111839173

$ws.Range("G57").Value = 0.1414749376738706
$ws.Range("H57").Value = 1.425833978070754
